$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.204.59'
$ws.Range('E2').Value = '  -3.71%  '
$ws.Range('D3').Value = '2.218.24'
$ws.Range('E3').Value = '  -6.70%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '296.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '82.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.510'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.47%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.467'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.62%  '
$ws.Range('E10').Value = '  -7.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '29.12'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.46'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -9.78%  '
$ws.Range('E13').Value = '  -2.10%  '
$ws.Range('D14').Value = '2.564.32'
$ws.Range('E14').Value = '  -6.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.22'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.12'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.30%  '
$ws.Range('D17').Value = '2.227.75'
$ws.Range('E17').Value = '  -7.24%  '
$ws.Range('E18').Value = '  -6.06%  '
$ws.Range('D19').Value = '39.103.85'
$ws.Range('E19').Value = '  -3.74%  '
$ws.Range('D20').Value = '0.0₃0870'
$ws.Range('E20').Value = '  -4.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.72'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.19'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '226.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.47%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('E26').Value = '  -7.23%  '
$ws.Range('E27').Value = '  -2.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.55'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.90%  '
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.09'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.48'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.84'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.43%  '
$ws.Range('E35').Value = '  -4.52%  '
$ws.Range('E36').Value = '  -4.53%  '
$ws.Range('E37').Value = '  -3.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.65'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.47%  '
$ws.Range('E39').Value = '  -4.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '14.95'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.12%  '
$ws.Range('E41').Value = '  -4.87%  '
$ws.Range('E42').Value = '  -3.42%  '
$ws.Range('D43').Value = '1.914.11'
$ws.Range('E43').Value = '  -2.62%  '
$ws.Range('E44').Value = '  -4.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.01'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -15.09%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.65'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.24%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.94'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.80'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -11.27%  '
$ws.Range('D49').Value = '2.435.87'
$ws.Range('E49').Value = '  -6.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '70.20'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '87.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.68%  '
